$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Pre-format the "cnt" column so new numeric-looking values are stored as text
# (matches the existing t="str" cell type used throughout column C).
$ws.Range("C2:C9").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = $true
$ws.Range("B2").Value = "edit"
$ws.Range("C2").Value = "3471"

# Row 3
$ws.Range("A3").Value = $false
$ws.Range("B3").Value = "edit"
$ws.Range("C3").Value = "2973"

# Row 4
$ws.Range("A4").Value = $false
$ws.Range("B4").Value = "categorize"
$ws.Range("C4").Value = "2161"

# Row 5
$ws.Range("A5").Value = $true
$ws.Range("B5").Value = "new"
$ws.Range("C5").Value = "747"

# Row 6
$ws.Range("A6").Value = $false
$ws.Range("B6").Value = "log"
$ws.Range("C6").Value = "355"

# Row 7
$ws.Range("A7").Value = $true
$ws.Range("B7").Value = "categorize"
$ws.Range("C7").Value = "228"

# Row 8 (newly inserted row)
$ws.Range("A8").Value = $true
$ws.Range("B8").Value = "log"
$ws.Range("C8").Value = "151"

# Row 9 (previously row 8, shifted down with an updated count)
$ws.Range("A9").Value = $false
$ws.Range("B9").Value = "new"
$ws.Range("C9").Value = "113"

Write-Host "done"
